$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 13.36918167694849
$ws.Cells.Item(2, 3).Value = 9.272512497947954
$ws.Cells.Item(2, 4).Value = 5.778062008106065
$ws.Cells.Item(2, 5).Value = 12.18380765609242
$ws.Cells.Item(2, 6).Value = 28.04071463641513
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 11).Value = 9.725229834376696
$ws.Cells.Item(2, 12).Value = 9.654118617904462
$ws.Cells.Item(2, 14).Value = 19.70236070378844
$ws.Cells.Item(2, 15).Value = 25.14028807524863
$ws.Cells.Item(3, 2).Value = 13.1311903008053
$ws.Cells.Item(3, 3).Value = 9.272639888474925
$ws.Cells.Item(3, 4).Value = 5.735678529978738
$ws.Cells.Item(3, 5).Value = 12.1935674632387
$ws.Cells.Item(3, 6).Value = 28.06402801069369
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 11).Value = 9.555774066044878
$ws.Cells.Item(3, 12).Value = 9.640379588668608
$ws.Cells.Item(3, 14).Value = 19.76338015506624
$ws.Cells.Item(3, 15).Value = 25.20186646751761
$ws.Cells.Item(4, 2).Value = 12.98559683549389
$ws.Cells.Item(4, 3).Value = 9.272997837982645
$ws.Cells.Item(4, 4).Value = 5.709035498388679
$ws.Cells.Item(4, 5).Value = 12.20178365642163
$ws.Cells.Item(4, 6).Value = 28.0850317534008
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 11).Value = 9.45191128257154
$ws.Cells.Item(4, 12).Value = 9.633608499674482
$ws.Cells.Item(4, 14).Value = 19.80259224988303
$ws.Cells.Item(4, 15).Value = 25.24461826394896
$ws.Cells.Item(5, 2).Value = 12.926481950924
$ws.Cells.Item(5, 3).Value = 9.273214445022559
$ws.Cells.Item(5, 4).Value = 5.698026013061111
$ws.Cells.Item(5, 5).Value = 12.20569123904671
$ws.Cells.Item(5, 6).Value = 28.09527168203759
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 11).Value = 9.409688178555736
$ws.Cells.Item(5, 12).Value = 9.631270221080007
$ws.Cells.Item(5, 14).Value = 19.81901180687404
$ws.Cells.Item(5, 15).Value = 25.263280901852
$ws.Cells.Item(6, 2).Value = 12.91668138473359
$ws.Cells.Item(6, 3).Value = 9.273254696618253
$ws.Cells.Item(6, 4).Value = 5.696188805371515
$ws.Cells.Item(6, 5).Value = 12.20637388196945
$ws.Cells.Item(6, 6).Value = 28.09707347280728
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 11).Value = 9.40268485173481
$ws.Cells.Item(6, 12).Value = 9.630907443025547
$ws.Cells.Item(6, 14).Value = 19.82176489333639
$ws.Cells.Item(6, 15).Value = 25.26645471008661
$ws.Cells.Item(7, 2).Value = 12.98479861037415
$ws.Cells.Item(7, 3).Value = 9.273000472270228
$ws.Cells.Item(7, 4).Value = 5.708887632342401
$ws.Cells.Item(7, 5).Value = 12.20183409016027
$ws.Cells.Item(7, 6).Value = 28.08516304976982
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 11).Value = 9.451341360798533
$ws.Cells.Item(7, 12).Value = 9.633575257290465
$ws.Cells.Item(7, 14).Value = 19.80281190526864
$ws.Cells.Item(7, 15).Value = 25.24486493292682
$ws.Cells.Item(8, 2).Value = 13.28705985826045
$ws.Cells.Item(8, 3).Value = 9.27249860278096
$ws.Cells.Item(8, 4).Value = 5.763578839534495
$ws.Cells.Item(8, 5).Value = 12.18671153180498
$ws.Cells.Item(8, 6).Value = 28.04736415673979
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 11).Value = 9.666796492801948
$ws.Cells.Item(8, 12).Value = 9.649037396978908
$ws.Cells.Item(8, 14).Value = 19.72303862707961
$ws.Cells.Item(8, 15).Value = 25.16049342345863
$ws.Cells.Item(9, 2).Value = 13.880519939901
$ws.Cells.Item(9, 3).Value = 9.273716838910687
$ws.Cells.Item(9, 4).Value = 5.865762446356865
$ws.Cells.Item(9, 5).Value = 12.17468212701295
$ws.Cells.Item(9, 6).Value = 28.02634610395992
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 11).Value = 10.08836208218175
$ws.Cells.Item(9, 12).Value = 9.692451669269319
$ws.Cells.Item(9, 14).Value = 19.58039658848075
$ws.Cells.Item(9, 15).Value = 25.03433577628521
$ws.Cells.Item(10, 2).Value = 14.31240521158942
$ws.Cells.Item(10, 3).Value = 9.275932052549807
$ws.Cells.Item(10, 4).Value = 5.937551829611305
$ws.Cells.Item(10, 5).Value = 12.17655675785786
$ws.Cells.Item(10, 6).Value = 28.04327429014531
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 11).Value = 10.39439118618921
$ws.Cells.Item(10, 12).Value = 9.732158788330493
$ws.Cells.Item(10, 14).Value = 19.48392081074432
$ws.Cells.Item(10, 15).Value = 24.96570272352362
$ws.Cells.Item(11, 2).Value = 14.50706629594722
$ws.Cells.Item(11, 3).Value = 9.277221694384387
$ws.Cells.Item(11, 4).Value = 5.969456681710195
$ws.Cells.Item(11, 5).Value = 12.17972496434753
$ws.Cells.Item(11, 6).Value = 28.05798632658014
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 11).Value = 10.53218297331897
$ws.Cells.Item(11, 12).Value = 9.751876467942697
$ws.Cells.Item(11, 14).Value = 19.44182059080072
$ws.Cells.Item(11, 15).Value = 24.93972091132089
$ws.Cells.Item(12, 2).Value = 14.58044749519691
$ws.Cells.Item(12, 3).Value = 9.277750173795177
$ws.Cells.Item(12, 4).Value = 5.981426595320507
$ws.Cells.Item(12, 5).Value = 12.18125636423936
$ws.Cells.Item(12, 6).Value = 28.06456257728833
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 11).Value = 10.58410726971051
$ws.Cells.Item(12, 12).Value = 9.759576965336951
$ws.Cells.Item(12, 14).Value = 19.42613397200689
$ws.Cells.Item(12, 15).Value = 24.93063685954002
$ws.Cells.Item(13, 2).Value = 14.56465946377114
$ws.Cells.Item(13, 3).Value = 9.27763457868064
$ws.Cells.Item(13, 4).Value = 5.978853686726811
$ws.Cells.Item(13, 5).Value = 12.18091182146943
$ws.Cells.Item(13, 6).Value = 28.06310161969428
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 11).Value = 10.57293652942041
$ws.Cells.Item(13, 12).Value = 9.757908193489921
$ws.Cells.Item(13, 14).Value = 19.42950100671261
$ws.Cells.Item(13, 15).Value = 24.93255968993847
$ws.Cells.Item(14, 2).Value = 14.51311048532632
$ws.Cells.Item(14, 3).Value = 9.277264371058703
$ws.Cells.Item(14, 4).Value = 5.970443713684607
$ws.Cells.Item(14, 5).Value = 12.17984431449438
$ws.Cells.Item(14, 6).Value = 28.05850723542898
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 11).Value = 10.5364601934563
$ws.Cells.Item(14, 12).Value = 9.752505330438266
$ws.Cells.Item(14, 14).Value = 19.44052492281626
$ws.Cells.Item(14, 15).Value = 24.93895842936876
$ws.Cells.Item(15, 2).Value = 14.48148980993376
$ws.Cells.Item(15, 3).Value = 9.277042820474753
$ws.Cells.Item(15, 4).Value = 5.965277703421769
$ws.Cells.Item(15, 5).Value = 12.17923358705984
$ws.Cells.Item(15, 6).Value = 28.05582383192406
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 11).Value = 10.51408280609755
$ws.Cells.Item(15, 12).Value = 9.749226247384305
$ws.Cells.Item(15, 14).Value = 19.44731067087346
$ws.Cells.Item(15, 15).Value = 24.94297615874865
$ws.Cells.Item(16, 2).Value = 14.29964098801476
$ws.Cells.Item(16, 3).Value = 9.275853408628288
$ws.Cells.Item(16, 4).Value = 5.935451325052694
$ws.Cells.Item(16, 5).Value = 12.17639619508255
$ws.Cells.Item(16, 6).Value = 28.04245366503954
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 11).Value = 10.38535322340602
$ws.Cells.Item(16, 12).Value = 9.730903143688197
$ws.Cells.Item(16, 14).Value = 19.48670802053652
$ws.Cells.Item(16, 15).Value = 24.96750624263552
$ws.Cells.Item(17, 2).Value = 14.18756701828559
$ws.Cells.Item(17, 3).Value = 9.275195653114238
$ws.Cells.Item(17, 4).Value = 5.916958526178137
$ws.Cells.Item(17, 5).Value = 12.17524773255671
$ws.Cells.Item(17, 6).Value = 28.03604534447157
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 11).Value = 10.30598119744867
$ws.Cells.Item(17, 12).Value = 9.720083487309337
$ws.Cells.Item(17, 14).Value = 19.51133390676475
$ws.Cells.Item(17, 15).Value = 24.98389761409528
$ws.Cells.Item(18, 2).Value = 14.12293922564206
$ws.Cells.Item(18, 3).Value = 9.274843891465274
$ws.Cells.Item(18, 4).Value = 5.906251337122556
$ws.Cells.Item(18, 5).Value = 12.17480519617138
$ws.Cells.Item(18, 6).Value = 28.03301962783156
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 11).Value = 10.26019746299532
$ws.Cells.Item(18, 12).Value = 9.714016319126012
$ws.Cells.Item(18, 14).Value = 19.52566634377319
$ws.Cells.Item(18, 15).Value = 24.99381862299775
$ws.Cells.Item(19, 2).Value = 14.10103117293323
$ws.Cells.Item(19, 3).Value = 9.274729365844275
$ws.Cells.Item(19, 4).Value = 5.902614047137591
$ws.Cells.Item(19, 5).Value = 12.17469284349917
$ws.Cells.Item(19, 6).Value = 28.03210863862299
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 11).Value = 10.24467493114993
$ws.Cells.Item(19, 12).Value = 9.711988996210035
$ws.Cells.Item(19, 14).Value = 19.53054799826023
$ws.Cells.Item(19, 15).Value = 24.99726235991065
$ws.Cells.Item(20, 2).Value = 14.19951518514556
$ws.Cells.Item(20, 3).Value = 9.275262926296604
$ws.Cells.Item(20, 4).Value = 5.918934451063003
$ws.Cells.Item(20, 5).Value = 12.17534743186174
$ws.Cells.Item(20, 6).Value = 28.03665921461712
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 11).Value = 10.31444441711304
$ws.Cells.Item(20, 12).Value = 9.721219140453602
$ws.Cells.Item(20, 14).Value = 19.50869503137081
$ws.Cells.Item(20, 15).Value = 24.98210167835214
$ws.Cells.Item(21, 2).Value = 14.52826126427368
$ws.Cells.Item(21, 3).Value = 9.277372024533225
$ws.Cells.Item(21, 4).Value = 5.972916987781954
$ws.Cells.Item(21, 5).Value = 12.18014887703736
$ws.Cells.Item(21, 6).Value = 28.05982946780627
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 11).Value = 10.54718146924045
$ws.Cells.Item(21, 12).Value = 9.754085969945086
$ws.Cells.Item(21, 14).Value = 19.43727999773177
$ws.Cells.Item(21, 15).Value = 24.9370584734762
$ws.Cells.Item(22, 2).Value = 14.74114515872961
$ws.Cells.Item(22, 3).Value = 9.278984095371429
$ws.Cells.Item(22, 4).Value = 6.007544665485994
$ws.Cells.Item(22, 5).Value = 12.18521937814674
$ws.Cells.Item(22, 6).Value = 28.08082912571273
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 11).Value = 10.69778399166001
$ws.Cells.Item(22, 12).Value = 9.776927255287397
$ws.Cells.Item(22, 14).Value = 19.39209673296483
$ws.Cells.Item(22, 15).Value = 24.91201950397996
$ws.Cells.Item(23, 2).Value = 14.62772852228639
$ws.Cells.Item(23, 3).Value = 9.278102463288661
$ws.Cells.Item(23, 4).Value = 5.989124130665966
$ws.Cells.Item(23, 5).Value = 12.18233681204727
$ws.Cells.Item(23, 6).Value = 28.06908658935566
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 11).Value = 10.6175581123449
$ws.Cells.Item(23, 12).Value = 9.764613331122526
$ws.Cells.Item(23, 14).Value = 19.41607588376593
$ws.Cells.Item(23, 15).Value = 24.92498037905641
$ws.Cells.Item(24, 2).Value = 14.19411402060211
$ws.Cells.Item(24, 3).Value = 9.27523242985718
$ws.Cells.Item(24, 4).Value = 5.91804137010477
$ws.Cells.Item(24, 5).Value = 12.17530167951209
$ws.Cells.Item(24, 6).Value = 28.03637963247343
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 11).Value = 10.31061866418095
$ws.Cells.Item(24, 12).Value = 9.720705234173591
$ws.Cells.Item(24, 14).Value = 19.50988752269398
$ws.Cells.Item(24, 15).Value = 24.98291207164164
$ws.Cells.Item(25, 2).Value = 13.72037684053851
$ws.Cells.Item(25, 3).Value = 9.273153803777568
$ws.Cells.Item(25, 4).Value = 5.838681877528209
$ws.Cells.Item(25, 5).Value = 12.17605202975851
$ws.Cells.Item(25, 6).Value = 28.02634389303911
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 11).Value = 9.974746828252266
$ws.Cells.Item(25, 12).Value = 9.679321820675799
$ws.Cells.Item(25, 14).Value = 19.61751725049802
$ws.Cells.Item(25, 15).Value = 25.06424696071539
